# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed values, as captured by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.160.13"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "2.324.89"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'521.63"
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").Value = "'135.08"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "2.348.90"
$ws.Range("E10").Value = "  +5.92%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "'5.32"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "2.742.50"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "56.980.81"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "2.327.44"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "'10.49"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "'323.31"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").Value = "'6.60"
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'60.82"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  +9.08%  "
$ws.Range("D26").Value = "'0.994"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'7.99"
$ws.Range("E27").Value = "  +6.17%  "
$ws.Range("E28").Value = "  +14.14%  "
$ws.Range("D29").Value = "0.0₃0745"
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("D30").Value = "'1.72"
$ws.Range("E30").Value = "  +5.32%  "
$ws.Range("D31").Value = "'168.16"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").Value = "'6.21"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").Value = "'0.930"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'4.03"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("E39").Value = "  +7.61%  "
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("D43").Value = "'138.62"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").Value = "'5.20"
$ws.Range("E44").Value = "  +5.90%  "
$ws.Range("D45").Value = "'277.42"
$ws.Range("E45").Value = "  +6.39%  "
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("D50").Value = "'17.85"
$ws.Range("E50").Value = "  +8.14%  "
$ws.Range("E51").Value = "  +0.81%  "
